# Raw pair angles service integration: add a new row (A4/B4) to Sheet1
# that mirrors the existing "rawOutagesCreationServiceUrl" hyperlink row,
# pointing at the same google.com placeholder URL, then select B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label + URL text for row 4 (adds a 6th shared string).
$ws.Range("A4").Value = "rawPairAnglesCreationServiceUrl"
$ws.Range("B4").Value = "http://google.com"

# Create the hyperlink relationship for B4 (mirrors B3 -> http://google.com/).
$ws.Hyperlinks.Add($ws.Range("B4"), "http://google.com/")

# Hyperlinks.Add() re-stamps B4's text/style; restore the plain display
# text + the shared "Hyperlink" cell style (same style B3 already uses).
$ws.Range("B4").Value = "http://google.com"
$ws.Range("B4").Style = "Hyperlink"

# Match the author's final selection state (cell B4 active/selected).
$ws.Range("B4").Select()
